$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (Steel) was missing its Immunities entry - fill in D18 with "Poison",
# reusing the existing shared string used elsewhere in column D.
$ws.Range("D18").Value = "Poison"

# Leave the selection on the cell that was just edited (D18), matching the
# saved workbook state.
$ws.Range("D18").Select()
